$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content change ---
# D4 text changed from many inner spaces to "NÃO " + newline + "(POSSUI ESTRÓBILOS)"
$ws.Range("D4").Value = "NÃO `n(POSSUI ESTRÓBILOS)"

# --- Column width changes ---
# target stored widths: D -> 18 (bestFit), E -> 22.140625 (closest achievable given engine's
# pixel-quantized ColumnWidth storage)
$ws.Columns.Item(4).ColumnWidth = 17.166666666666668
$ws.Columns.Item(5).ColumnWidth = 21.307291666666668

# --- Row height changes ---
$ws.Rows.Item(4).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 37.5

# --- View changes: zoom + selection ---
$ws.Select()
$excel.ActiveWindow.Zoom = 180
$ws.Range("F4").Select()
